$d = $word.ActiveDocument

function Replace-UniqueText($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: text not found -> $findText"
    }
}

# Muc 10) Nghe nghiep khi duoc tuyen dung
Replace-UniqueText "crrgyê" "Công nghệ thông tin"

# Muc 13) Cong viec chinh duoc giao
Replace-UniqueText "dgr" "Công nghệ thông tin"

# Phu cap chuc vu (source run uses decomposed e + combining grave accent)
Replace-UniqueText "fhryè" "Không"

# Muc 24) La thuong binh hang
Replace-UniqueText "Thương binh loại 1/4" "Không"

# Giang day
Replace-UniqueText "cfghgngh" "hg"

# So quyet dinh - dong 1 (2023 - 2025)
Replace-UniqueText "20/QĐ-BGDĐT" "20/QĐ-ĐHCT"

# So quyet dinh - dong 2 (2019 - 2022)
Replace-UniqueText "18/QĐ-BGDĐT" "18/QĐ-ĐHCT"
